$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value reads as a plain number (e.g. "211.50", "0.512").
# Force a Text number format before writing so Excel keeps the literal
# string (with trailing zeros / precision) instead of coercing it to a
# numeric value, then clear the formatting so no stray style lingers.
$numericLooking = @(
    'D5', 'D6', 'D10', 'D14', 'D16', 'D20', 'D21', 'D22',
    'D25', 'D27', 'D29', 'D32', 'D35', 'D36', 'D39', 'D40',
    'D43', 'D44', 'D46', 'D47', 'D49', 'D50', 'D51'
)
foreach ($ref in $numericLooking) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply the updated coin snapshot (prices + 1h volume deltas, and for a
# few rows the coin name/link where the ranking order shifted).
$ws.Range('D2').Value = '26.685.01'
$ws.Range('E2').Value = '  +0.00%  '
$ws.Range('D3').Value = '1.596.57'
$ws.Range('E3').Value = '  +0.02%  '
$ws.Range('D5').Value = '211.50'
$ws.Range('E5').Value = '  +0.26%  '
$ws.Range('D6').Value = '0.512'
$ws.Range('E6').Value = '  -0.20%  '
$ws.Range('E7').Value = '  +0.27%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E9').Value = '  +0.60%  '
$ws.Range('D10').Value = '19.40'
$ws.Range('E10').Value = '  -0.56%  '
$ws.Range('E11').Value = '  -0.43%  '
$ws.Range('D12').Value = '1.820.65'
$ws.Range('E12').Value = '  +0.02%  '
$ws.Range('D13').Value = '1.624.95'
$ws.Range('E13').Value = '  +1.85%  '
$ws.Range('D14').Value = '4.02'
$ws.Range('E14').Value = '  +0.10%  '
$ws.Range('E15').Value = '  +0.44%  '
$ws.Range('D16').Value = '65.00'
$ws.Range('E16').Value = '  +0.50%  '
$ws.Range('D17').Value = '26.653.02'
$ws.Range('E17').Value = '  -0.09%  '
$ws.Range('D18').Value = '0.0₃0751'
$ws.Range('E18').Value = '  +3.12%  '
$ws.Range('E19').Value = '  +0.28%  '
$ws.Range('D20').Value = '208.95'
$ws.Range('E20').Value = '  +0.59%  '
$ws.Range('D21').Value = '6.96'
$ws.Range('E21').Value = '  +2.68%  '
$ws.Range('D22').Value = '4.26'
$ws.Range('E22').Value = '  +0.57%  '
$ws.Range('E23').Value = '  -1.14%  '
$ws.Range('E24').Value = '  +1.08%  '
$ws.Range('D25').Value = '142.93'
$ws.Range('E25').Value = '  -1.75%  '
$ws.Range('E26').Value = '  +0.23%  '
$ws.Range('D27').Value = '7.11'
$ws.Range('E27').Value = '  -1.30%  '
$ws.Range('E28').Value = '  -0.64%  '
$ws.Range('D29').Value = '15.31'
$ws.Range('E29').Value = '  +0.55%  '
$ws.Range('E30').Value = '  +2.49%  '
$ws.Range('E31').Value = '  +0.07%  '
$ws.Range('D32').Value = '3.23'
$ws.Range('E32').Value = '  +0.30%  '
$ws.Range('E33').Value = '  +0.59%  '
$ws.Range('D34').Value = '1.281.62'
$ws.Range('E34').Value = '  -0.16%  '
$ws.Range('D35').Value = '0.614'
$ws.Range('E35').Value = '  -7.33%  '
$ws.Range('D36').Value = '2.45'
$ws.Range('E36').Value = '  -0.39%  '
$ws.Range('E37').Value = '  -0.38%  '
$ws.Range('E38').Value = '  +0.04%  '
$ws.Range('B39').Value = 'WEMIXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D39').Value = '1.05'
$ws.Range('E39').Value = '  +16.59%  '
$ws.Range('B40').Value = 'ARBITRUM'
$ws.Range('C40').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D40').Value = '0.824'
$ws.Range('E40').Value = '  -2.06%  '
$ws.Range('E41').Value = '  +0.13%  '
$ws.Range('E42').Value = '  -0.33%  '
$ws.Range('D43').Value = '0.782'
$ws.Range('E43').Value = '  -0.43%  '
$ws.Range('D44').Value = '62.74'
$ws.Range('E44').Value = '  -1.11%  '
$ws.Range('D45').Value = '1.733.54'
$ws.Range('E45').Value = '  +0.00%  '
$ws.Range('D46').Value = '90.87'
$ws.Range('E46').Value = '  +0.96%  '
$ws.Range('D47').Value = '1.56'
$ws.Range('E47').Value = '  -2.77%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').Value = '0.0₆0104'
$ws.Range('E48').Value = '  -1.38%  '
$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').Value = '0.100'
$ws.Range('E49').Value = '  -0.54%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = '0.0509'
$ws.Range('E50').Value = '  +0.72%  '
$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D51').Value = '1.00'
$ws.Range('E51').Value = '  +0.24%  '

foreach ($ref in $numericLooking) {
    $ws.Range($ref).ClearFormats()
}
